# Update coded-segment timestamps in column M: normalize 2-digit year
# (e.g. "6/18/18 00:07:00") to 4-digit year (e.g. "6/18/2018 00:07:00").
# These are plain text cells (number format "Text"), so we assign the
# corrected text directly -- no date re-parsing should occur.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "6/18/2018 00:07:00"
$ws.Range("M3").Value = "6/18/2018 00:08:00"
$ws.Range("M4").Value = "6/18/2018 00:36:00"
$ws.Range("M5").Value = "6/18/2018 00:36:00"
$ws.Range("M6").Value = "6/18/2018 00:36:00"
$ws.Range("M7").Value = "6/18/2018 00:37:00"
$ws.Range("M8").Value = "6/18/2018 00:37:00"
$ws.Range("M9").Value = "6/18/2018 01:01:00"
$ws.Range("M10").Value = "6/18/2018 01:02:00"
$ws.Range("M11").Value = "6/18/2018 01:02:00"
$ws.Range("M12").Value = "6/18/2018 01:15:00"
$ws.Range("M13").Value = "6/18/2018 01:15:00"
$ws.Range("M14").Value = "6/18/2018 01:15:00"
$ws.Range("M15").Value = "6/18/2018 01:16:00"
$ws.Range("M16").Value = "6/18/2018 01:16:00"
$ws.Range("M17").Value = "6/18/2018 01:16:00"
$ws.Range("M18").Value = "6/18/2018 01:17:00"
$ws.Range("M19").Value = "6/18/2018 01:17:00"
$ws.Range("M20").Value = "6/18/2018 01:17:00"
$ws.Range("M21").Value = "6/18/2018 01:17:00"
$ws.Range("M22").Value = "6/18/2018 01:17:00"
$ws.Range("M23").Value = "6/18/2018 01:18:00"
$ws.Range("M24").Value = "6/18/2018 01:18:00"
$ws.Range("M25").Value = "6/18/2018 01:23:00"
$ws.Range("M26").Value = "6/18/2018 01:23:00"
$ws.Range("M27").Value = "6/18/2018 01:24:00"
$ws.Range("M28").Value = "6/18/2018 01:24:00"
$ws.Range("M29").Value = "6/18/2018 01:24:00"
$ws.Range("M30").Value = "6/18/2018 01:24:00"
$ws.Range("M31").Value = "6/18/2018 01:30:00"
$ws.Range("M32").Value = "6/18/2018 01:30:00"
$ws.Range("M33").Value = "6/18/2018 01:30:00"
$ws.Range("M34").Value = "6/18/2018 01:30:00"
$ws.Range("M35").Value = "6/18/2018 01:30:00"
$ws.Range("M36").Value = "6/18/2018 01:30:00"
$ws.Range("M37").Value = "6/18/2018 01:31:00"
$ws.Range("M38").Value = "6/18/2018 01:31:00"
$ws.Range("M39").Value = "6/18/2018 01:36:00"
$ws.Range("M40").Value = "6/18/2018 01:36:00"
$ws.Range("M41").Value = "6/18/2018 01:36:00"
$ws.Range("M42").Value = "6/18/2018 01:36:00"
$ws.Range("M43").Value = "6/18/2018 01:36:00"
$ws.Range("M44").Value = "6/18/2018 01:37:00"
$ws.Range("M45").Value = "6/18/2018 01:37:00"
$ws.Range("M46").Value = "6/18/2018 01:37:00"
$ws.Range("M47").Value = "6/18/2018 01:37:00"
$ws.Range("M48").Value = "6/18/2018 01:37:00"
$ws.Range("M49").Value = "6/18/2018 01:37:00"
$ws.Range("M50").Value = "6/18/2018 01:37:00"
$ws.Range("M51").Value = "6/18/2018 01:37:00"
$ws.Range("M52").Value = "6/18/2018 01:37:00"
$ws.Range("M53").Value = "6/18/2018 01:38:00"
$ws.Range("M54").Value = "6/18/2018 01:39:00"
$ws.Range("M55").Value = "6/18/2018 01:39:00"
$ws.Range("M56").Value = "6/18/2018 01:39:00"
$ws.Range("M57").Value = "6/18/2018 01:39:00"
$ws.Range("M58").Value = "6/18/2018 01:40:00"
$ws.Range("M59").Value = "6/18/2018 01:40:00"
$ws.Range("M60").Value = "6/18/2018 01:40:00"
$ws.Range("M61").Value = "6/18/2018 01:41:00"
$ws.Range("M62").Value = "6/18/2018 01:41:00"
$ws.Range("M63").Value = "6/18/2018 01:42:00"
$ws.Range("M64").Value = "6/18/2018 01:42:00"
$ws.Range("M65").Value = "6/18/2018 09:18:00"
$ws.Range("M66").Value = "6/18/2018 09:18:00"
$ws.Range("M67").Value = "6/18/2018 09:18:00"
$ws.Range("M68").Value = "6/18/2018 09:19:00"
$ws.Range("M69").Value = "6/18/2018 09:19:00"
$ws.Range("M70").Value = "6/18/2018 09:21:00"
$ws.Range("M71").Value = "6/18/2018 09:22:00"
$ws.Range("M72").Value = "6/18/2018 09:22:00"
$ws.Range("M73").Value = "6/18/2018 09:23:00"
$ws.Range("M74").Value = "6/18/2018 09:23:00"
$ws.Range("M75").Value = "6/18/2018 09:24:00"
$ws.Range("M76").Value = "6/18/2018 09:24:00"
$ws.Range("M77").Value = "6/18/2018 09:24:00"
$ws.Range("M78").Value = "6/18/2018 09:24:00"
$ws.Range("M79").Value = "6/18/2018 09:24:00"
$ws.Range("M80").Value = "6/18/2018 09:24:00"
$ws.Range("M81").Value = "6/18/2018 09:25:00"
$ws.Range("M82").Value = "6/18/2018 09:25:00"
$ws.Range("M83").Value = "6/18/2018 09:25:00"
$ws.Range("M84").Value = "6/18/2018 09:27:00"
$ws.Range("M85").Value = "6/18/2018 09:28:00"
$ws.Range("M86").Value = "6/18/2018 09:34:00"
$ws.Range("M87").Value = "6/18/2018 09:34:00"
$ws.Range("M88").Value = "6/18/2018 09:35:00"
$ws.Range("M89").Value = "6/18/2018 09:35:00"
$ws.Range("M90").Value = "6/18/2018 09:35:00"
$ws.Range("M91").Value = "6/18/2018 09:35:00"
$ws.Range("M92").Value = "6/18/2018 09:35:00"
$ws.Range("M93").Value = "6/18/2018 09:43:00"
$ws.Range("M94").Value = "6/18/2018 09:44:00"
$ws.Range("M95").Value = "6/18/2018 09:44:00"
$ws.Range("M96").Value = "6/18/2018 09:45:00"
$ws.Range("M97").Value = "6/18/2018 09:45:00"
$ws.Range("M98").Value = "6/18/2018 09:53:00"
$ws.Range("M99").Value = "6/18/2018 09:55:00"
$ws.Range("M100").Value = "6/18/2018 09:55:00"
$ws.Range("M101").Value = "6/18/2018 09:55:00"
$ws.Range("M102").Value = "6/18/2018 09:55:00"
$ws.Range("M103").Value = "6/18/2018 09:55:00"
$ws.Range("M104").Value = "6/18/2018 09:55:00"
$ws.Range("M105").Value = "6/18/2018 09:55:00"
$ws.Range("M106").Value = "6/18/2018 09:56:00"
$ws.Range("M107").Value = "6/18/2018 09:56:00"
$ws.Range("M108").Value = "6/18/2018 09:56:00"
$ws.Range("M109").Value = "6/18/2018 09:56:00"
$ws.Range("M110").Value = "6/18/2018 09:56:00"
$ws.Range("M111").Value = "6/18/2018 09:56:00"
$ws.Range("M112").Value = "6/18/2018 09:56:00"
$ws.Range("M113").Value = "6/18/2018 09:56:00"
$ws.Range("M114").Value = "6/18/2018 09:57:00"
$ws.Range("M115").Value = "6/18/2018 11:23:00"
$ws.Range("M116").Value = "6/18/2018 11:24:00"
$ws.Range("M117").Value = "6/18/2018 11:24:00"
$ws.Range("M118").Value = "6/18/2018 11:35:00"
$ws.Range("M119").Value = "6/18/2018 11:36:00"
$ws.Range("M120").Value = "10/29/2018 11:59:00"
$ws.Range("M121").Value = "10/29/2018 12:26:00"
$ws.Range("M122").Value = "10/30/2018 15:55:00"
$ws.Range("M123").Value = "10/30/2018 15:55:00"
$ws.Range("M124").Value = "10/30/2018 15:55:00"
$ws.Range("M125").Value = "10/30/2018 15:55:00"
$ws.Range("M126").Value = "10/30/2018 15:55:00"
$ws.Range("M127").Value = "10/30/2018 15:56:00"
$ws.Range("M128").Value = "10/30/2018 15:56:00"
$ws.Range("M129").Value = "10/30/2018 15:56:00"
$ws.Range("M130").Value = "10/30/2018 15:56:00"
$ws.Range("M131").Value = "10/30/2018 15:56:00"
$ws.Range("M132").Value = "10/30/2018 15:56:00"
$ws.Range("M133").Value = "10/30/2018 15:56:00"
$ws.Range("M134").Value = "10/30/2018 15:56:00"
$ws.Range("M135").Value = "10/30/2018 15:56:00"
$ws.Range("M136").Value = "10/30/2018 15:57:00"
$ws.Range("M137").Value = "10/30/2018 15:57:00"
$ws.Range("M138").Value = "10/30/2018 15:57:00"
$ws.Range("M139").Value = "10/30/2018 15:57:00"
$ws.Range("M140").Value = "10/30/2018 15:57:00"
$ws.Range("M141").Value = "10/30/2018 15:57:00"
$ws.Range("M142").Value = "10/30/2018 15:57:00"
$ws.Range("M143").Value = "10/30/2018 15:57:00"
$ws.Range("M144").Value = "10/30/2018 15:57:00"
$ws.Range("M145").Value = "11/8/2018 14:04:00"
$ws.Range("M146").Value = "11/8/2018 14:05:00"
$ws.Range("M147").Value = "11/12/2018 12:31:00"
$ws.Range("M148").Value = "11/12/2018 12:32:00"
$ws.Range("M149").Value = "11/12/2018 12:32:00"
$ws.Range("M150").Value = "1/29/2019 16:38:51"
$ws.Range("M151").Value = "8/22/2019 14:16:18"
$ws.Range("M152").Value = "8/22/2019 14:16:22"
$ws.Range("M153").Value = "8/22/2019 14:16:28"
